# Tidy up tests: clear stale "Targetable" (y) markers in the Parameters
# sheet's E column (rows 2-6), which are leftover/incorrect values now
# that the column header is "Targetable" but these parameters aren't
# actually flagged that way, and update the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parameters")
$ws.Activate()

# Clear the stale values in E2:E6 (previously held the shared string "y")
$ws.Range("E2:E6").ClearContents()

# Update the selected range to match the new view state
$ws.Range("E2:E8").Select()
